# Auto update Excel log
# Appends newly logged sensor readings (2026-02-06) to the PIR, Humidity and
# Temperature sheets, extending each sheet's used range.

$wb = $excel.ActiveWorkbook

function Append-SensorRows {
    param(
        [string]$SheetName,
        [int]$StartRow,
        [array]$Rows,
        [bool]$ForceTextValue
    )

    $ws = $wb.Worksheets.Item($SheetName)
    $endRow = $StartRow + $Rows.Count - 1

    # Column A holds plain-text dates like "2026-02-06"; without forcing a
    # text number format Excel will silently convert them into date serials.
    $ws.Range("A" + $StartRow + ":A" + $endRow).NumberFormat = "@"

    if ($ForceTextValue) {
        # Column E on this sheet holds percentage-looking text ("67.6%")
        # that Excel would otherwise coerce into a numeric percentage.
        $ws.Range("E" + $StartRow + ":E" + $endRow).NumberFormat = "@"
    }

    $r = $StartRow
    foreach ($row in $Rows) {
        $ws.Range("A" + $r).Value = "2026-02-06"
        $ws.Range("B" + $r).Value = $row[0]
        $ws.Range("C" + $r).Value = $row[1]
        $ws.Range("D" + $r).Value = $row[2]
        $ws.Range("E" + $r).Value = $row[3]
        $ws.Range("F" + $r).Value = $row[4]
        $r = $r + 1
    }
}

# ---- PIR sheet: rows 580-593 (motion / "No Motion" events) ----
$pirRows = @(
    @("10:28:08","10:00","Bathroom","No Motion","Inactive"),
    @("10:28:12","10:00","Bathroom","No Motion","Inactive"),
    @("10:28:16","10:00","Bathroom","No Motion","Inactive"),
    @("10:28:20","10:00","Bathroom","No Motion","Inactive"),
    @("10:28:23","10:00","Bathroom","No Motion","Inactive"),
    @("10:28:28","10:00","Bathroom","No Motion","Inactive"),
    @("10:28:33","10:00","Bathroom","No Motion","Inactive"),
    @("10:28:38","10:00","Bathroom","No Motion","Inactive"),
    @("10:28:43","10:00","Bathroom","No Motion","Inactive"),
    @("10:28:48","10:00","Bathroom","No Motion","Inactive"),
    @("10:28:53","10:00","Bathroom","No Motion","Inactive"),
    @("10:28:58","10:00","Bathroom","No Motion","Inactive"),
    @("10:29:03","10:00","Bathroom","No Motion","Inactive"),
    @("10:29:08","10:00","Bathroom","No Motion","Inactive")
)
Append-SensorRows "PIR" 580 $pirRows $false

# ---- Humidity sheet: rows 412-423 ----
$humidityRows = @(
    @("10:28:10","10:00","Bathroom","67.6%","Active"),
    @("10:28:14","10:00","Bathroom","67.6%","Active"),
    @("10:28:17","10:00","Bathroom","67.5%","Active"),
    @("10:28:24","10:00","Bathroom","67.5%","Active"),
    @("10:28:29","10:00","Bathroom","67.6%","Active"),
    @("10:28:34","10:00","Bathroom","67.7%","Active"),
    @("10:28:39","10:00","Bathroom","67.6%","Active"),
    @("10:28:44","10:00","Bathroom","67.7%","Active"),
    @("10:28:49","10:00","Bathroom","67.8%","Active"),
    @("10:28:54","10:00","Bathroom","67.8%","Active"),
    @("10:28:59","10:00","Bathroom","67.7%","Active"),
    @("10:29:04","10:00","Bathroom","67.7%","Active")
)
Append-SensorRows "Humidity" 412 $humidityRows $true

# ---- Temperature sheet: rows 412-423 ----
$temperatureRows = @(
    @("10:28:11","10:00","Bathroom","28.4C","Active"),
    @("10:28:15","10:00","Bathroom","28.5C","Active"),
    @("10:28:19","10:00","Bathroom","28.5C","Active"),
    @("10:28:25","10:00","Bathroom","28.4C","Active"),
    @("10:28:30","10:00","Bathroom","28.4C","Active"),
    @("10:28:35","10:00","Bathroom","28.5C","Active"),
    @("10:28:40","10:00","Bathroom","28.4C","Active"),
    @("10:28:45","10:00","Bathroom","28.4C","Active"),
    @("10:28:51","10:00","Bathroom","28.4C","Active"),
    @("10:28:55","10:00","Bathroom","28.4C","Active"),
    @("10:29:00","10:00","Bathroom","28.4C","Active"),
    @("10:29:06","10:00","Bathroom","28.4C","Active")
)
Append-SensorRows "Temperature" 412 $temperatureRows $false
